# Staging.Output.xlsx: header row was regenerated by the staging-template
# tool after the file moved under StagingTemplates/ -- the column headers
# in row 2 are now alphabetically sorted, "OutcomeSourceKey" was renamed to
# "OutcomeBusinessKey", and a few cosmetic/workbook-identity fields changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 header values: re-sorted alphabetically, OutcomeSourceKey -> OutcomeBusinessKey ---
$ws.Range("A2").Value = "BusinessKey"
$ws.Range("B2").Value = "Code"
$ws.Range("C2").Value = "LongName"
$ws.Range("D2").Value = "OutcomeBusinessKey"
$ws.Range("E2").Value = "Output_ID"
$ws.Range("F2").Value = "ShortName"
$ws.Range("G2").Value = "TextDescription"

# --- Cosmetic workbook-identity fields touched by the same save in Excel ---
# (best-effort: harmless no-ops if the host doesn't surface these via COM)
try { $ws.CodeName = "Sheet38" } catch {}
try {
    $win = $wb.Windows.Item(1)
    $win.Width = 28800
    $win.Height = 12585
} catch {}
